# "Add all entities db"
# The former "dia_fijo" column header (C1) is replaced by a new
# "Tipo de materia" column header; the other headers/values are untouched.
# Once "dia_fijo" is no longer referenced anywhere in the workbook it drops
# out of the shared-string table and "Tipo de materia" is appended at the
# end of that table on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Tipo de materia"

# Widen column C (new header is long) and nudge columns B/D to their new
# saved widths.
$ws.Columns.Item(2).ColumnWidth = 13.5
$ws.Columns.Item(3).ColumnWidth = 13.166666666666666
$ws.Columns.Item(4).ColumnWidth = 21

# Move the active selection to A7, matching the saved view state.
[void]$ws.Range("A7").Select()
